# Penalty Reward System (unfinished) — shift the 16-week forecast window
# forward by one week (drop W1 2025-01-05, add a new W16 2025-04-27) and
# refresh the dependent MyForecast values + Summary roll-ups.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Sheet "Forecast Comparison": Week_Start_Date (col B) is stored as
# plain text in the source file, not a real date. Format the column as
# text FIRST so Excel doesn't silently coerce the "YYYY-MM-DD" strings
# into date serials when we assign them.
$ws1.Range("B2:B17").NumberFormat = "@"

$ws1.Cells.Item(2,2).Value  = "2025-01-12"
$ws1.Cells.Item(2,4).Value  = 25

$ws1.Cells.Item(3,2).Value  = "2025-01-19"
$ws1.Cells.Item(3,4).Value  = 26

$ws1.Cells.Item(4,2).Value  = "2025-01-26"
$ws1.Cells.Item(4,4).Value  = 27

$ws1.Cells.Item(5,2).Value  = "2025-02-02"
$ws1.Cells.Item(5,4).Value  = 27

$ws1.Cells.Item(6,2).Value  = "2025-02-09"
$ws1.Cells.Item(6,4).Value  = 27

$ws1.Cells.Item(7,2).Value  = "2025-02-16"
$ws1.Cells.Item(7,4).Value  = 28

$ws1.Cells.Item(8,2).Value  = "2025-02-23"
$ws1.Cells.Item(8,4).Value  = 29

$ws1.Cells.Item(9,2).Value  = "2025-03-02"
$ws1.Cells.Item(9,4).Value  = 29

$ws1.Cells.Item(10,2).Value = "2025-03-09"
$ws1.Cells.Item(10,4).Value = 30

$ws1.Cells.Item(11,2).Value = "2025-03-16"
$ws1.Cells.Item(11,4).Value = 30

$ws1.Cells.Item(12,2).Value = "2025-03-23"
$ws1.Cells.Item(12,4).Value = 31

$ws1.Cells.Item(13,2).Value = "2025-03-30"
$ws1.Cells.Item(13,4).Value = 32

$ws1.Cells.Item(14,2).Value = "2025-04-06"
$ws1.Cells.Item(14,4).Value = 32

$ws1.Cells.Item(15,2).Value = "2025-04-13"
$ws1.Cells.Item(15,4).Value = 32

$ws1.Cells.Item(16,2).Value = "2025-04-20"
$ws1.Cells.Item(16,4).Value = 33

$ws1.Cells.Item(17,2).Value = "2025-04-27"
$ws1.Cells.Item(17,4).Value = 34

# --- Sheet "Summary": every cell in column B is text in the source file
# (including ones that look like plain numbers, e.g. "14"), so force text
# formatting first to keep that representation intact. Only touch the
# rows that actually change so untouched rows keep their original style.
# (Set per-cell rather than a comma-joined multi-area range — multi-area
# NumberFormat assignment only reliably hits the first area here.)
foreach ($r in @(2,5,6,8,9,10,11,12,13,14,15)) {
    $ws2.Cells.Item($r,2).NumberFormat = "@"
}

$ws2.Range("B2").Value  = "2024-05-12 to 2025-01-05"
$ws2.Range("B5").Value  = "13"
$ws2.Range("B6").Value  = "13"
$ws2.Range("B8").Value  = "472 units"
$ws2.Range("B9").Value  = "472"
$ws2.Range("B10").Value = "218"
$ws2.Range("B11").Value = "105"
$ws2.Range("B12").Value = "34"
$ws2.Range("B13").Value = "2025-04-27"
$ws2.Range("B14").Value = "25"
$ws2.Range("B15").Value = "2025-01-12"
